$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.4
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 6.5
$ws.Range("R2").Value = 1.91
$ws.Range("S2").Value = 1.91
$ws.Range("T2").Value = 8.5
$ws.Range("U2").Value = 7.5
$ws.Range("W2").Value = 10
$ws.Range("Z2").Value = 15
$ws.Range("AB2").Value = 21
$ws.Range("AD2").Value = 401
$ws.Range("AE2").Value = 19
$ws.Range("AG2").Value = 21
$ws.Range("AH2").Value = 81
$ws.Range("AI2").Value = 51
$ws.Range("G3").Value = 2.25
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 3.1
$ws.Range("L3").Value = 1.25
$ws.Range("M3").Value = 4
$ws.Range("T3").Value = 10
$ws.Range("X3").Value = 19
$ws.Range("AB3").Value = 13
$ws.Range("AD3").Value = 201
$ws.Range("AH3").Value = 34
$ws.Range("G4").Value = 1.18
$ws.Range("H4").Value = 8.5
$ws.Range("J4").Value = 1.01
$ws.Range("K4").Value = 23
$ws.Range("L4").Value = 1.11
$ws.Range("M4").Value = 6.5
$ws.Range("N4").Value = 1.36
$ws.Range("O4").Value = 3.2
$ws.Range("P4").Value = 1.2
$ws.Range("Q4").Value = 4.33
$ws.Range("W4").Value = 8
$ws.Range("Z4").Value = 23
$ws.Range("AA4").Value = 17
$ws.Range("AC4").Value = 81
$ws.Range("AD4").Value = 451
$ws.Range("G5").Value = 1.73
$ws.Range("H5").Value = 4.2
$ws.Range("I5").Value = 4.2
$ws.Range("N5").Value = 1.44
$ws.Range("O5").Value = 2.75
$ws.Range("AA5").Value = 9.5
$ws.Range("AB5").Value = 13
$ws.Range("AE5").Value = 21
$ws.Range("AI5").Value = 34
$ws.Range("G6").Value = 6
$ws.Range("I6").Value = 1.55
$ws.Range("J6").Value = 1.05
$ws.Range("K6").Value = 11
$ws.Range("N6").Value = 1.86
$ws.Range("O6").Value = 2.04
$ws.Range("P6").Value = 1.36
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 1.91
$ws.Range("S6").Value = 1.91
$ws.Range("V6").Value = 21
$ws.Range("X6").Value = 51
$ws.Range("Y6").Value = 51
$ws.Range("AD6").Value = 351
$ws.Range("AF6").Value = 8
$ws.Range("AH6").Value = 12
$ws.Range("AI6").Value = 15
$ws.Range("AJ6").Value = 29
$ws.Range("G7").Value = 1.5
$ws.Range("I7").Value = 6.25
$ws.Range("N7").Value = 1.95
$ws.Range("O7").Value = 1.95
$ws.Range("G8").Value = 1.44
$ws.Range("J8").Value = 1.04
$ws.Range("K8").Value = 13
$ws.Range("L8").Value = 1.25
$ws.Range("M8").Value = 3.75
$ws.Range("N8").Value = 1.88
$ws.Range("O8").Value = 1.98
$ws.Range("K26").Value = 23
$ws.Range("L26").Value = 1.11
$ws.Range("M26").Value = 6.5
$ws.Range("N26").Value = 1.4
$ws.Range("O26").Value = 2.88
$ws.Range("P26").Value = 1.22
$ws.Range("Q26").Value = 4
$ws.Range("R26").Value = 2.05
$ws.Range("S26").Value = 1.7
$ws.Range("U26").Value = 7
$ws.Range("V26").Value = 11
$ws.Range("W26").Value = 7
$ws.Range("Y26").Value = 29
$ws.Range("Z26").Value = 21
$ws.Range("AD26").Value = 401
$ws.Range("G27").Value = 3.6
$ws.Range("I27").Value = 1.91
$ws.Range("L27").Value = 1.4
$ws.Range("M27").Value = 2.75
$ws.Range("T27").Value = 9.5
$ws.Range("V27").Value = 13
$ws.Range("AH27").Value = 17
$ws.Range("AI27").Value = 19
$ws.Range("H28").Value = 3.75
$ws.Range("I28").Value = 2.63
$ws.Range("J28").Value = 1.03
$ws.Range("K28").Value = 10
$ws.Range("N28").Value = 1.7
$ws.Range("O28").Value = 2.1
$ws.Range("P28").Value = 1.33
$ws.Range("Q28").Value = 3.25
$ws.Range("T28").Value = 9.5
$ws.Range("X28").Value = 17
$ws.Range("AA28").Value = 7.5
$ws.Range("AG28").Value = 10
$ws.Range("AJ28").Value = 26
$ws.Range("G29").Value = 1.57
$ws.Range("H29").Value = 4.33
$ws.Range("J29").Value = 1.03
$ws.Range("K29").Value = 15
$ws.Range("N29").Value = 1.73
$ws.Range("O29").Value = 2.08
$ws.Range("X29").Value = 12
$ws.Range("AF29").Value = 29
$ws.Range("AG29").Value = 17
$ws.Range("G30").Value = 2.1
$ws.Range("H30").Value = 3.5
$ws.Range("I30").Value = 3.3
$ws.Range("J30").Value = 1.05
$ws.Range("K30").Value = 11
$ws.Range("N30").Value = 1.93
$ws.Range("O30").Value = 1.93
$ws.Range("P30").Value = 1.36
$ws.Range("Q30").Value = 3
$ws.Range("R30").Value = 1.67
$ws.Range("S30").Value = 2.1
$ws.Range("T30").Value = 8.5
$ws.Range("U30").Value = 11
$ws.Range("Y30").Value = 23
$ws.Range("Z30").Value = 12
$ws.Range("AA30").Value = 7
$ws.Range("AB30").Value = 13
$ws.Range("AC30").Value = 41
$ws.Range("AD30").Value = 151
$ws.Range("AE30").Value = 11
$ws.Range("AJ30").Value = 29
